# EggDescriptor.xlsx update
# - Rewrites the "Internal Functionality" example rows (previously the
#   "Eg Turn Left / Eg Push / ShouldTurnLeft / Return" placeholder text)
#   with the real EggControl script description (Remain Idle / NotInIce /
#   Time rows), and sets the InterFace name to "EggControl".
# - Widens columns B and C so the new, longer descriptions are readable.
# - Leaves the final selection on B16 (where the author's edit ended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: was the lone "Eg Turn Left" example; now the start of the real
# "Remain Idle..." description, with a matching Text Description in D5.
$ws.Range("A5").Value = "Remain Idle at set location or inside a"
$ws.Range("D5").Value = "waits for game conditions to be met to trigger spawn enemy script"

# New rows 6-7 continue the multi-line description started in A5.
$ws.Range("A6").Value = "parented ice cube until conditions like "
$ws.Range("A7").Value = "time or ice cube destroyed are met"

# Row 10: External Outgoing example - "Eg Push" replaced by the real
# "NotInIce" condition, plus its Text Description and Parameters columns.
$ws.Range("A10").Value = "NotInIce"
$ws.Range("B10").Value = "NotInIce: Confirms the egg is not in an iceblock"
$ws.Range("C10").Value = "spawnEnemy(IceDestroyed = true && NotInIce = true || time > 3 seconds && NotInIce = true)"

# New row 11: second External Outgoing entry, "Time".
$ws.Range("A11").Value = "Time "
$ws.Range("B11").Value = "Time: Hatch egg after a few seconds"

# Row 14: External Incoming example - "ShouldTurnLeft"/"Return" replaced
# with the real "IceDestroyed" condition and its description.
$ws.Range("A14").Value = "IceDestroyed "
$ws.Range("C14").Value = "Confirms an iceblock has been destroyed "

# Row 16: Name of InterFace value updated from the old tutorial link to
# the actual interface name, "EggControl".
$ws.Range("B16").Value = "EggControl"

# Widen columns B and C to fit the longer text now held in them
# (target stored widths: B=41.140625, C=80.140625 "characters"; the
# ColumnWidth setter snaps to whole-pixel increments, so we feed in the
# values that land on the closest achievable pixel width).
$ws.Columns("B").ColumnWidth = 40.3
$ws.Columns("C").ColumnWidth = 79.33

# Leave the selection on B16, matching the end state of the author's edit.
$ws.Range("B16").Select()
